$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.9179281773574478
$ws.Range("J2").Value = 0.9179281773574478
$ws.Range("M2").Value = 4.407279333333333
$ws.Range("N2").Value = 13.221838
$ws.Range("O2").Value = 0.2946616623342344
$ws.Range("P2").Value = 0.2946616623342344
$ws.Range("Q2").Value = 2.718311463561556
$ws.Range("R2").Value = 24.464803172054
$ws.Range("S2").Value = 0.2704782426435794
$ws.Range("T2").Value = 0.2704782426435794
$ws.Range("I3").Value = 0.9179281773574478
$ws.Range("J3").Value = 0.9179281773574478
$ws.Range("O3").Value = 0.2393683991842171
$ws.Range("P3").Value = 0.2393683991842171
$ws.Range("S3").Value = 0.2197229983801384
$ws.Range("T3").Value = 0.2197229983801384
$ws.Range("I4").Value = 0.9179281773574478
$ws.Range("J4").Value = 0.9179281773574478
$ws.Range("M4").Value = 3.580339
$ws.Range("N4").Value = 10.741017
$ws.Range("O4").Value = 0.2393741266819538
$ws.Range("P4").Value = 0.2393741266819538
$ws.Range("Q4").Value = 2.208273134295667
$ws.Range("R4").Value = 19.874458208661
$ws.Range("S4").Value = 0.2197282558116967
$ws.Range("T4").Value = 0.2197282558116967
$ws.Range("I5").Value = 0.9179281773574478
$ws.Range("J5").Value = 0.9179281773574478
$ws.Range("M5").Value = 3.389212666666667
$ws.Range("N5").Value = 10.167638
$ws.Range("O5").Value = 0.2265958117995947
$ws.Range("P5").Value = 0.2265958117995947
$ws.Range("Q5").Value = 2.090390680383778
$ws.Range("R5").Value = 18.813516123454
$ws.Range("S5").Value = 0.2079986805220333
$ws.Range("T5").Value = 0.2079986805220333
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.055146
$ws.Range("H6").Value = 0.165438
$ws.Range("I6").Value = 0.08207182264255215
$ws.Range("J6").Value = 0.08207182264255215
$ws.Range("M6").Value = 4.407279333333333
$ws.Range("N6").Value = 13.221838
$ws.Range("O6").Value = 0.2946616623342344
$ws.Range("P6").Value = 0.2946616623342344
$ws.Range("Q6").Value = 0.243043826116
$ws.Range("R6").Value = 2.187394435044
$ws.Range("S6").Value = 0.02418341969065487
$ws.Range("T6").Value = 0.02418341969065487
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.055146
$ws.Range("H7").Value = 0.165438
$ws.Range("I7").Value = 0.08207182264255215
$ws.Range("J7").Value = 0.08207182264255215
$ws.Range("O7").Value = 0.2393683991842171
$ws.Range("P7").Value = 0.2393683991842171
$ws.Range("Q7").Value = 0.19743665032
$ws.Range("R7").Value = 1.77692985288
$ws.Range("S7").Value = 0.01964540080407869
$ws.Range("T7").Value = 0.01964540080407869
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.055146
$ws.Range("H8").Value = 0.165438
$ws.Range("I8").Value = 0.08207182264255215
$ws.Range("J8").Value = 0.08207182264255215
$ws.Range("M8").Value = 3.580339
$ws.Range("N8").Value = 10.741017
$ws.Range("O8").Value = 0.2393741266819538
$ws.Range("P8").Value = 0.2393741266819538
$ws.Range("Q8").Value = 0.197441374494
$ws.Range("R8").Value = 1.776972370446
$ws.Range("S8").Value = 0.01964587087025712
$ws.Range("T8").Value = 0.01964587087025713
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.055146
$ws.Range("H9").Value = 0.165438
$ws.Range("I9").Value = 0.08207182264255215
$ws.Range("J9").Value = 0.08207182264255215
$ws.Range("M9").Value = 3.389212666666667
$ws.Range("N9").Value = 10.167638
$ws.Range("O9").Value = 0.2265958117995947
$ws.Range("P9").Value = 0.2265958117995947
$ws.Range("Q9").Value = 0.186901521716
$ws.Range("R9").Value = 1.682113695444
$ws.Range("S9").Value = 0.01859713127756147
$ws.Range("T9").Value = 0.01859713127756147
